# Dp.xlsx edit: add "Coin Exchange" worksheet after "LIS", with the
# coin-exchange DP table values, a thin box border style, column widths,
# and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the last sheet (after "LIS") ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Coin Exchange"

# --- Column widths (A:H narrow, J wide) -----------------------------------
$ws.Columns("A:H").ColumnWidth = 1.8
$ws.Columns("J:J").ColumnWidth = 12

# --- Row 2: coin weights header row, boxed with a thin border ------------
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 2

$headerRng = $ws.Range("B2:H2")
$headerRng.Borders.LineStyle = 1
$headerRng.Borders.Weight = 2

# --- Row 3: running money total (0..6) ------------------------------------
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 6

# --- Problem statement labels ---------------------------------------------
$ws.Range("J2").Value = "coins:1,3,4"
$ws.Range("K2").Value = "money:6"

# --- Page setup (matches the other sheets) --------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / active cell ----------------------------------------------
[void]$ws.Range("K3").Select()

# --- Make the new sheet the active tab (this also clears tabSelected on
#     the previously active "LIS" sheet) -----------------------------------
[void]$ws.Activate()
